# Error Calculations and Plots
# Removes two sample rows ("RM 232" and "SC 92") from the missing-data
# sheet, which shifts all subsequent rows up by two, and refreshes a
# handful of individual B/C/F measurement cells (some newly populated,
# some cleared back to missing) to match the updated imputation pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that no longer belong in this sample set ---
# Row 26 is "RM 232"; deleting it shifts "SC 92" (originally row 28) up
# to row 27, so the second delete targets row 27 to remove "SC 92".
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Cleared cells (value removed / now missing) ---
$ws.Range("F5").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("B26").ClearContents()
$ws.Range("C27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("F29").ClearContents()

# --- Newly populated / corrected cells ---
$ws.Range("F11").Value = 17.65
$ws.Range("C19").Value = 13.2
$ws.Range("C23").Value = 12.2
$ws.Range("F23").Value = 16.48
$ws.Range("F25").Value = 16.6
$ws.Range("B27").Value = -20.4
$ws.Range("C33").Value = 10.4
$ws.Range("F33").Value = 17.53
